$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ticker rows that are no longer part of the trimmed example list.
# (Deleted bottom-up so earlier row numbers stay valid.)
$ws.Range("F22").EntireRow.Delete()   | Out-Null   # GPS / THE GAP
$ws.Range("F18:F20").EntireRow.Delete() | Out-Null # JWN/NORDSTROM, PVH/..., RL/RALPH LAUREN
$ws.Range("F13:F15").EntireRow.Delete() | Out-Null # FL/FOOT LOCKER, GIL/GILDAN ACTIVEWEAR, LB/L BRANDS
$ws.Range("F11").EntireRow.Delete()   | Out-Null   # AEO / AMERICAN EAGLE OUTFITTERS

# Resulting table is now: TICKER/COMPANY NAME header, BURL, M, NKE/NIKE, ROST.

# Narrow column G and drop its "best fit" autosize flag.
$ws.Columns("G").ColumnWidth = 21.3

# Match the saved selection from the edit.
$ws.Range("G4").Select() | Out-Null

Write-Output "done"
